$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New experiment rows (10-14) - fill Serial No / Experiment / Date of Performance
# Row 12 -> Serial 10
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Memory Allocation Methods"
$ws.Range("C12").Value = 45420

# Row 13 -> Serial 11
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Page Replacement Algorithms"
$ws.Range("C13").Value = 45421

# Row 14 -> Serial 12
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Deadlock Avoidance - Banker's Algorithm"
$ws.Range("C14").Value = 45421

# Row 15 -> Serial 13
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Deadlock Detection Algorithm"
$ws.Range("C15").Value = 45421

# Row 16 -> Serial 14
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Disk Scheduling algorithm"
$ws.Range("C16").Value = 45421

# The date cells need the same number-format/style as the existing date column
# (re-apply formatting from an already-formatted date cell so the stored
# cell style index is reused rather than a brand new style being created).
$ws.Range("C3").Copy()
$ws.Range("C12:C16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the saved selection to match the author's final cursor position
$ws.Range("B5").Select()
